# Generate Report for Handoff
# Adds two newly-discovered source files (279d1811-... and e5c85d09-...) to the
# localization-status report, in between the existing e8ee0bed-... row and the
# 8e99a18b-... row (alphabetical-ish ordering), and appends the e5c85d09-... row
# right after 8e99a18b-..., before the trailing .localization-config row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Hyperlinks don't automatically track row shifts caused by Rows.Insert(), so
# drop them all up front and recreate them at the end once every row is in
# its final location.
$ws1.Hyperlinks.Delete()

# Before: row2=e8ee0bed, row3=8e99a18b, row4=.localization-config
# Insert a row at 3 (pushes 8e99a18b -> row4, .localization-config -> row5)
$ws1.Rows.Item(3).Insert()
# Insert a row at 5 (pushes .localization-config -> row6)
$ws1.Rows.Item(5).Insert()

$ws1.Range("A3").Value = "279d1811-555e-490b-ac6e-300ca84faa46.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A5").Value = "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.md"
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/29d3a927d18928fc149d8440d50a37b932613fbb/e2e/e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md", "", "", "e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5e2375c68c0bfe2189e2bd653a5deefdca7eacf3/e2e/279d1811-555e-490b-ac6e-300ca84faa46.md", "", "", "279d1811-555e-490b-ac6e-300ca84faa46.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a1acef856f2938ef635038727c5638883b5e6d48/e2e/8e99a18b-4369-4c7c-92e2-73849d6401c9.md", "", "", "8e99a18b-4369-4c7c-92e2-73849d6401c9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0241be560beacf215cdfdf2adcd958ebe6c8b008/e2e/e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.md", "", "", "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/29d3a927d18928fc149d8440d50a37b932613fbb/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Rows.Item(3).Insert()
$ws2.Rows.Item(5).Insert()

$ws2.Range("A3").Value = "279d1811-555e-490b-ac6e-300ca84faa46.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "279d1811-555e-490b-ac6e-300ca84faa46.5e2375c68c0bfe2189e2bd653a5deefdca7eacf3.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-10 06:53:33"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A5").Value = "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.md"
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.0241be560beacf215cdfdf2adcd958ebe6c8b008.zh-cn.xlf"
$ws2.Range("D5").Value = "0001-01-01 00:00:00"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/29d3a927d18928fc149d8440d50a37b932613fbb/e2e/e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md", "", "", "e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0d1315bdb17350eb2e1b1722809ae19a557c8006/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.65c4accda60b6940698eb1ab94c4a6ac937ad02d.zh-cn.xlf", "", "", "e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.65c4accda60b6940698eb1ab94c4a6ac937ad02d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6e656f7a08f0fa689b5474104e04d5c2b800727e/e2e/e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md", "", "", "e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/993c95fcdb9404e37c5b3e6fdf5cc2213bc6007e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.65c4accda60b6940698eb1ab94c4a6ac937ad02d.zh-cn.xlf", "", "", "e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.65c4accda60b6940698eb1ab94c4a6ac937ad02d.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5e2375c68c0bfe2189e2bd653a5deefdca7eacf3/e2e/279d1811-555e-490b-ac6e-300ca84faa46.md", "", "", "279d1811-555e-490b-ac6e-300ca84faa46.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e2375c68c0bfe2189e2bd653a5deefdca7eacf3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/279d1811-555e-490b-ac6e-300ca84faa46.5e2375c68c0bfe2189e2bd653a5deefdca7eacf3.zh-cn.xlf", "", "", "279d1811-555e-490b-ac6e-300ca84faa46.5e2375c68c0bfe2189e2bd653a5deefdca7eacf3.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a1acef856f2938ef635038727c5638883b5e6d48/e2e/8e99a18b-4369-4c7c-92e2-73849d6401c9.md", "", "", "8e99a18b-4369-4c7c-92e2-73849d6401c9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4f13d15c12e4e802d448c5d6827a4b02c9c5b13d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8e99a18b-4369-4c7c-92e2-73849d6401c9.eb99b62922e0e19cca5f70210ab6149c0c813899.zh-cn.xlf", "", "", "8e99a18b-4369-4c7c-92e2-73849d6401c9.eb99b62922e0e19cca5f70210ab6149c0c813899.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0241be560beacf215cdfdf2adcd958ebe6c8b008/e2e/e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.md", "", "", "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0241be560beacf215cdfdf2adcd958ebe6c8b008/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.0241be560beacf215cdfdf2adcd958ebe6c8b008.zh-cn.xlf", "", "", "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.0241be560beacf215cdfdf2adcd958ebe6c8b008.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/29d3a927d18928fc149d8440d50a37b932613fbb/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Rows.Item(3).Insert()
$ws3.Rows.Item(5).Insert()

$ws3.Range("A3").Value = "279d1811-555e-490b-ac6e-300ca84faa46.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "279d1811-555e-490b-ac6e-300ca84faa46.5e2375c68c0bfe2189e2bd653a5deefdca7eacf3.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-10 06:53:39"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A5").Value = "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.md"
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.0241be560beacf215cdfdf2adcd958ebe6c8b008.de-de.xlf"
$ws3.Range("D5").Value = "2016-03-10 06:53:39"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/29d3a927d18928fc149d8440d50a37b932613fbb/e2e/e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md", "", "", "e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7dd8d70240a6258c064dfbf3bad814413388ffff/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.65c4accda60b6940698eb1ab94c4a6ac937ad02d.de-de.xlf", "", "", "e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.65c4accda60b6940698eb1ab94c4a6ac937ad02d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/cfd22f3070c5185e50e6b4b4d3eef98a26afa1dd/e2e/e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md", "", "", "e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/be9641dc7463eee8e80d92413710a650e3e03c94/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.65c4accda60b6940698eb1ab94c4a6ac937ad02d.de-de.xlf", "", "", "e8ee0bed-07cf-4ded-9ce3-32b10dfc202c.65c4accda60b6940698eb1ab94c4a6ac937ad02d.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5e2375c68c0bfe2189e2bd653a5deefdca7eacf3/e2e/279d1811-555e-490b-ac6e-300ca84faa46.md", "", "", "279d1811-555e-490b-ac6e-300ca84faa46.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e2375c68c0bfe2189e2bd653a5deefdca7eacf3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/279d1811-555e-490b-ac6e-300ca84faa46.5e2375c68c0bfe2189e2bd653a5deefdca7eacf3.de-de.xlf", "", "", "279d1811-555e-490b-ac6e-300ca84faa46.5e2375c68c0bfe2189e2bd653a5deefdca7eacf3.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a1acef856f2938ef635038727c5638883b5e6d48/e2e/8e99a18b-4369-4c7c-92e2-73849d6401c9.md", "", "", "8e99a18b-4369-4c7c-92e2-73849d6401c9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d6c3943777117bff455fba191a16859effc3aff6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8e99a18b-4369-4c7c-92e2-73849d6401c9.eb99b62922e0e19cca5f70210ab6149c0c813899.de-de.xlf", "", "", "8e99a18b-4369-4c7c-92e2-73849d6401c9.eb99b62922e0e19cca5f70210ab6149c0c813899.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0241be560beacf215cdfdf2adcd958ebe6c8b008/e2e/e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.md", "", "", "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0241be560beacf215cdfdf2adcd958ebe6c8b008/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.0241be560beacf215cdfdf2adcd958ebe6c8b008.de-de.xlf", "", "", "e5c85d09-ceac-4d7a-b5b0-bfad6a442bb0.0241be560beacf215cdfdf2adcd958ebe6c8b008.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/29d3a927d18928fc149d8440d50a37b932613fbb/.localization-config", "", "", ".localization-config") | Out-Null
